# =============================================================
# CompStat weekly report refresh: new crime data collected
# =============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: issue number and week-covering date range ---
$ws.Range("A8").Value() = "Volume 32   Number  14"
$ws.Range("C9").Value() = "Report Covering the Week  3/31/2025  Through  4/6/2025"

# --- Column H width grows to match column E (both now show 3-digit pct values) ---
$ws.Columns.Item(8).ColumnWidth() = 6.65

# --- Cells that become the placeholder text ("0" / "***.*") ---
# Reference cells C14 (-> "0", shared string 20) and E14 (-> "***.*", shared string 21)
# already carry the desired text + style; copy their format, then their value, onto the target.
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$excel.CutCopyMode() = 0

# --- Updated crime-statistic figures ---
# Row 15
$ws.Range("F15").Value() = 3
$ws.Range("H15").Value() = 200
# Row 16
$ws.Range("D16").Value() = 5
$ws.Range("E16").Value() = -100
$ws.Range("G16").Value() = 9
$ws.Range("H16").Value() = 0
$ws.Range("I16").Value() = 26
$ws.Range("J16").Value() = 29
$ws.Range("K16").Value() = -10.344827586206
$ws.Range("L16").Value() = -25.714285714285
$ws.Range("M16").Value() = 44.444444444444
$ws.Range("N16").Value() = -83.647798742138
# Row 17
$ws.Range("C17").Value() = 2
$ws.Range("D17").Value() = 4
$ws.Range("E17").Value() = -50
$ws.Range("F17").Value() = 10
$ws.Range("G17").Value() = 9
$ws.Range("H17").Value() = 11.111111111111
$ws.Range("I17").Value() = 34
$ws.Range("J17").Value() = 25
$ws.Range("K17").Value() = 36
$ws.Range("L17").Value() = 3.030303030303
$ws.Range("M17").Value() = -19.047619047619
$ws.Range("N17").Value() = -47.692307692307
# Row 18
$ws.Range("F18").Value() = 8
$ws.Range("G18").Value() = 7
$ws.Range("H18").Value() = 14.285714285714
$ws.Range("L18").Value() = 3.846153846153
$ws.Range("M18").Value() = -12.903225806451
$ws.Range("N18").Value() = -78.571428571428
# Row 19
$ws.Range("C19").Value() = 12
$ws.Range("D19").Value() = 16
$ws.Range("F19").Value() = 50
$ws.Range("G19").Value() = 52
$ws.Range("H19").Value() = -3.846153846153
$ws.Range("I19").Value() = 162
$ws.Range("J19").Value() = 171
$ws.Range("K19").Value() = -5.263157894736
$ws.Range("L19").Value() = -12.432432432432
$ws.Range("M19").Value() = 5.194805194805
$ws.Range("N19").Value() = -27.027027027027
# Row 20
$ws.Range("D20").Value() = 1
$ws.Range("E20").Value() = -100
$ws.Range("G20").Value() = 4
$ws.Range("H20").Value() = -50
$ws.Range("J20").Value() = 15
$ws.Range("K20").Value() = -73.333333333333
$ws.Range("L20").Value() = -63.636363636363
$ws.Range("M20").Value() = -66.666666666666
$ws.Range("N20").Value() = -97.163120567375
# Row 21
$ws.Range("C21").Value() = 14
$ws.Range("D21").Value() = 26
$ws.Range("E21").Value() = -46.153846153846
$ws.Range("F21").Value() = 82
$ws.Range("G21").Value() = 82
$ws.Range("H21").Value() = 0
$ws.Range("I21").Value() = 259
$ws.Range("J21").Value() = 270
$ws.Range("K21").Value() = -4.074074074074
$ws.Range("L21").Value() = -11.301369863013
$ws.Range("M21").Value() = 0.387596899224
$ws.Range("N21").Value() = -64.027777777777
# Row 22
$ws.Range("C22").Value() = 1
$ws.Range("E22").Value() = 0
$ws.Range("F22").Value() = 3
$ws.Range("H22").Value() = 0
$ws.Range("I22").Value() = 9
$ws.Range("J22").Value() = 12
$ws.Range("K22").Value() = -25
$ws.Range("L22").Value() = -10
$ws.Range("M22").Value() = 200
# Row 23
$ws.Range("D23").Value() = 1
$ws.Range("E23").Value() = -100
$ws.Range("F23").Value() = 4
$ws.Range("H23").Value() = 33.333333333333
$ws.Range("J23").Value() = 11
$ws.Range("K23").Value() = 36.363636363636
$ws.Range("L23").Value() = 25
$ws.Range("M23").Value() = 7.142857142857
# Row 24
$ws.Range("C24").Value() = 23
$ws.Range("D24").Value() = 9
$ws.Range("E24").Value() = 155.555555555556
$ws.Range("F24").Value() = 73
$ws.Range("G24").Value() = 43
$ws.Range("H24").Value() = 69.767441860465
$ws.Range("I24").Value() = 212
$ws.Range("J24").Value() = 152
$ws.Range("K24").Value() = 39.473684210526
$ws.Range("L24").Value() = 21.142857142857
$ws.Range("M24").Value() = 6.532663316582
# Row 25
$ws.Range("C25").Value() = 11
$ws.Range("E25").Value() = 57.142857142857
$ws.Range("F25").Value() = 40
$ws.Range("G25").Value() = 17
$ws.Range("H25").Value() = 135.294117647059
$ws.Range("I25").Value() = 116
$ws.Range("J25").Value() = 67
$ws.Range("K25").Value() = 73.134328358209
$ws.Range("L25").Value() = 13.725490196078
# Row 26
$ws.Range("C26").Value() = 5
$ws.Range("D26").Value() = 6
$ws.Range("E26").Value() = -16.666666666666
$ws.Range("G26").Value() = 22
$ws.Range("H26").Value() = -4.545454545454
$ws.Range("I26").Value() = 67
$ws.Range("J26").Value() = 87
$ws.Range("K26").Value() = -22.988505747126
$ws.Range("L26").Value() = -28.723404255319
$ws.Range("M26").Value() = -16.25
# Row 27
$ws.Range("D27").Value() = 1
$ws.Range("E27").Value() = 0
$ws.Range("G27").Value() = 2
$ws.Range("H27").Value() = 100
$ws.Range("I27").Value() = 7
$ws.Range("J27").Value() = 7
$ws.Range("L27").Value() = 133.333333333333
# Row 28
$ws.Range("C28").Value() = 1
$ws.Range("G28").Value() = 8
$ws.Range("H28").Value() = -50
$ws.Range("I28").Value() = 12
$ws.Range("K28").Value() = -33.333333333333
$ws.Range("L28").Value() = -29.411764705882
# Row 29
$ws.Range("D29").Value() = 1
$ws.Range("E29").Value() = -100
$ws.Range("G29").Value() = 1
$ws.Range("H29").Value() = -100
$ws.Range("J29").Value() = 1
$ws.Range("K29").Value() = 0
# Row 30
$ws.Range("D30").Value() = 1
$ws.Range("E30").Value() = -100
$ws.Range("G30").Value() = 1
$ws.Range("H30").Value() = -100
$ws.Range("J30").Value() = 1
$ws.Range("K30").Value() = 0
